# Swap the presentation's theme color scheme from the "Integral" (Red
# Violet) palette over to the stock "Office Theme" palette.
#
# (The source deck ships two theme parts: theme1.xml - "Integral"/"Red
# Violet" - wired to the slide master that every slide actually uses,
# and theme2.xml - "Office Theme" - wired only to the notes master. The
# authored change swaps their contents outright. The notes-master theme
# isn't reachable as a distinct object through this PowerPoint host - it
# resolves back to the same Theme as the slide master - so we drive the
# visible/applied theme's ThemeColorScheme to the Office palette, which
# is the change that is actually observable on every slide.)

function ToComRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = ToComRGB($officeColors[$i - 1])
}
